# Rename header suffixes:
#   *_old -> *_FV2404
#   *_new -> *_FV2410
# then wrap the data range in an Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($oldHeaders[$i] -replace "_old$", "_FV2404")
}

# Column 11 is "diff" and stays unchanged; new headers begin at column 12 (L)
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($newHeaders[$i] -replace "_new$", "_FV2410")
}

$rng = $ws.Range("A1:U62")
$listObj = $ws.ListObjects.Add(1, $rng, 0, 1)
$listObj.Name = "Table1"

$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
